$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: new rows 4 and 5 ---
$overview.Range("A4").Value = "13119eb3-2e3c-4293-a369-ceaae48b7a54.md"
$overview.Range("B4").Value = 'e2e\13119eb3-2e3c-4293-a369-ceaae48b7a54.md'
$overview.Range("C4").Value = ".md"
$overview.Range("D4").Value = ""
$overview.Range("E4").Value = "Ready for handoff"
$overview.Range("F4").Value = "Ready for handoff"
$overview.Range("G4").Value = "2016-10-25 02:07:47"

$overview.Range("A5").Value = "ad829fe1-993c-47c3-a276-61654914ba8f.md"
$overview.Range("B5").Value = 'e2e\ad829fe1-993c-47c3-a276-61654914ba8f.md'
$overview.Range("C5").Value = ".md"
$overview.Range("D5").Value = ""
$overview.Range("E5").Value = "Ready for handoff"
$overview.Range("F5").Value = "Ready for handoff"
$overview.Range("G5").Value = "2016-10-25 02:07:47"

$overview.Hyperlinks.Add($overview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/13119eb3-2e3c-4293-a369-ceaae48b7a54.md", [Type]::Missing, [Type]::Missing, 'e2e\13119eb3-2e3c-4293-a369-ceaae48b7a54.md')
$overview.Hyperlinks.Add($overview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/ad829fe1-993c-47c3-a276-61654914ba8f.md", [Type]::Missing, [Type]::Missing, 'e2e\ad829fe1-993c-47c3-a276-61654914ba8f.md')

# --- zh-cn sheet: new rows 4 and 5 ---
$zhcn.Range("A4").Value = "13119eb3-2e3c-4293-a369-ceaae48b7a54.md"
$zhcn.Range("B4").Value = ".md"
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("D4").Value = "e2e"
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("F4").Value = "False"
$zhcn.Range("G4").Value = "13119eb3-2e3c-4293-a369-ceaae48b7a54.442b3eeac2579df5dea2580826eb0a62258b20a3.zh-cn.xlf"
$zhcn.Range("H4").Value = "2016-10-25 02:07:35"
$zhcn.Range("I4").Value = ""
$zhcn.Range("J4").Value = ""
$zhcn.Range("K4").Value = "0001-01-01 00:00:00"
$zhcn.Range("L4").Value = ""
$zhcn.Range("M4").Value = "True"
$zhcn.Range("N4").Value = ""
$zhcn.Range("O4").Value = "False"
$zhcn.Range("P4").Value = ""

$zhcn.Range("A5").Value = "ad829fe1-993c-47c3-a276-61654914ba8f.md"
$zhcn.Range("B5").Value = ".md"
$zhcn.Range("C5").Value = "Ready for handoff"
$zhcn.Range("D5").Value = "e2e"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("F5").Value = "False"
$zhcn.Range("G5").Value = "ad829fe1-993c-47c3-a276-61654914ba8f.379f31a8c9df08a612250e737b0b9a4be3f86d36.zh-cn.xlf"
$zhcn.Range("H5").Value = "2016-10-25 02:07:35"
$zhcn.Range("I5").Value = ""
$zhcn.Range("J5").Value = ""
$zhcn.Range("K5").Value = "0001-01-01 00:00:00"
$zhcn.Range("L5").Value = ""
$zhcn.Range("M5").Value = "True"
$zhcn.Range("N5").Value = ""
$zhcn.Range("O5").Value = "False"
$zhcn.Range("P5").Value = ""

$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/13119eb3-2e3c-4293-a369-ceaae48b7a54.md", [Type]::Missing, [Type]::Missing, "13119eb3-2e3c-4293-a369-ceaae48b7a54.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/ad829fe1-993c-47c3-a276-61654914ba8f.md", [Type]::Missing, [Type]::Missing, "ad829fe1-993c-47c3-a276-61654914ba8f.md")

# --- de-de sheet: new rows 4 and 5 ---
$dede.Range("A4").Value = "13119eb3-2e3c-4293-a369-ceaae48b7a54.md"
$dede.Range("B4").Value = ".md"
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("D4").Value = "e2e"
$dede.Range("E4").Value = "ht"
$dede.Range("F4").Value = "False"
$dede.Range("G4").Value = "13119eb3-2e3c-4293-a369-ceaae48b7a54.442b3eeac2579df5dea2580826eb0a62258b20a3.de-de.xlf"
$dede.Range("H4").Value = "2016-10-25 02:07:47"
$dede.Range("I4").Value = ""
$dede.Range("J4").Value = ""
$dede.Range("K4").Value = "0001-01-01 00:00:00"
$dede.Range("L4").Value = ""
$dede.Range("M4").Value = "True"
$dede.Range("N4").Value = ""
$dede.Range("O4").Value = "False"
$dede.Range("P4").Value = ""

$dede.Range("A5").Value = "ad829fe1-993c-47c3-a276-61654914ba8f.md"
$dede.Range("B5").Value = ".md"
$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("D5").Value = "e2e"
$dede.Range("E5").Value = "ht"
$dede.Range("F5").Value = "False"
$dede.Range("G5").Value = "ad829fe1-993c-47c3-a276-61654914ba8f.379f31a8c9df08a612250e737b0b9a4be3f86d36.de-de.xlf"
$dede.Range("H5").Value = "2016-10-25 02:07:47"
$dede.Range("I5").Value = ""
$dede.Range("J5").Value = ""
$dede.Range("K5").Value = "0001-01-01 00:00:00"
$dede.Range("L5").Value = ""
$dede.Range("M5").Value = "True"
$dede.Range("N5").Value = ""
$dede.Range("O5").Value = "False"
$dede.Range("P5").Value = ""

$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/13119eb3-2e3c-4293-a369-ceaae48b7a54.md", [Type]::Missing, [Type]::Missing, "13119eb3-2e3c-4293-a369-ceaae48b7a54.md")
$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2d1b692cfd814447656bd35c542b5ae7d1e75d7/e2e/ad829fe1-993c-47c3-a276-61654914ba8f.md", [Type]::Missing, [Type]::Missing, "ad829fe1-993c-47c3-a276-61654914ba8f.md")

# --- Resize the tables to cover the new rows ---
$overview.ListObjects.Item(1).Resize($overview.Range("A1:G5"))
$zhcn.ListObjects.Item(1).Resize($zhcn.Range("A1:P5"))
$dede.ListObjects.Item(1).Resize($dede.Range("A1:P5"))

# --- Column width changes ---
$overview.Range("E1").EntireColumn.ColumnWidth = 17.2159881591797
$overview.Range("F1").EntireColumn.ColumnWidth = 17.2159881591797
$zhcn.Range("C1").EntireColumn.ColumnWidth = 17.2159881591797
$dede.Range("C1").EntireColumn.ColumnWidth = 17.2159881591797

Write-Host "done"
